# Applies the "Updated symbol list" price/volume refresh described in the commit.
# Source workbook stores every data cell as literal text (inline strings), so
# numeric-looking values are written with a leading apostrophe to force Excel to
# keep them as Text instead of auto-converting to Number (which would lose
# trailing zeros / flip to scientific notation for the very small price values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'326.84"
$ws.Range("E2").Value = "'0.00%"

# Row 3
$ws.Range("E3").Value = "'-0.99%"

# Row 4
$ws.Range("D4").Value = "'5.507"
$ws.Range("E4").Value = "'-0.72%"

# Row 5
$ws.Range("D5").Value = "'0.08012"
$ws.Range("E5").Value = "'-0.91%"

# Row 6
$ws.Range("D6").Value = "'1.993"
$ws.Range("E6").Value = "'4.36%"

# Row 7
$ws.Range("E7").Value = "'-0.93%"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9529"
$ws.Range("E8").Value = "'0.43%"

# Row 9
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.568"
$ws.Range("E9").Value = "'-5.04%"

# Row 10
$ws.Range("D10").Value = "'0.1125"
$ws.Range("E10").Value = "'-4.47%"

# Row 11
$ws.Range("E11").Value = "'-1.75%"

# Row 12
$ws.Range("D12").Value = "'10.59"
$ws.Range("E12").Value = "'26.75%"

# Row 13
$ws.Range("D13").Value = "'0.09881"
$ws.Range("E13").Value = "'-1.99%"

# Row 14
$ws.Range("D14").Value = "'0.04593"
$ws.Range("E14").Value = "'9.83%"

# Row 15
$ws.Range("D15").Value = "'0.1066"
$ws.Range("E15").Value = "'0.10%"

# Row 16
$ws.Range("D16").Value = "'0.001256"
$ws.Range("E16").Value = "'-0.90%"

# Row 17
$ws.Range("D17").Value = "'0.04086"
$ws.Range("E17").Value = "'-3.91%"

# Row 18
$ws.Range("D18").Value = "'0.005966"
$ws.Range("E18").Value = "'0.38%"

# Row 19
$ws.Range("D19").Value = "'3.354"
$ws.Range("E19").Value = "'-6.74%"

# Row 21
$ws.Range("D21").Value = "'0.1406"
$ws.Range("E21").Value = "'2.32%"

# Row 22
$ws.Range("D22").Value = "'0.2544"
$ws.Range("E22").Value = "'-4.55%"

# Row 23
$ws.Range("D23").Value = "'0.001260"
$ws.Range("E23").Value = "'1.60%"

# Row 24
$ws.Range("D24").Value = "'0.004318"
$ws.Range("E24").Value = "'-6.04%"

# Row 25
$ws.Range("D25").Value = "'0.0001161"
$ws.Range("E25").Value = "'-5.79%"

# Row 26
$ws.Range("D26").Value = "'0.0003746"
$ws.Range("E26").Value = "'-6.38%"

# Row 38
$ws.Range("D38").Value = "'0.02548"
$ws.Range("E38").Value = "'-4.03%"

# Row 39
$ws.Range("D39").Value = "'0.05620"
$ws.Range("E39").Value = "'1.30%"

# Row 40
$ws.Range("D40").Value = "'0.007542"
$ws.Range("E40").Value = "'-1.92%"

# Row 41
$ws.Range("D41").Value = "'0.1397"
$ws.Range("E41").Value = "'0.22%"

# Row 42
$ws.Range("D42").Value = "'0.007595"
$ws.Range("E42").Value = "'-32.99%"

# Row 43
$ws.Range("D43").Value = "'0.002016"
$ws.Range("E43").Value = "'-2.03%"

# Row 44
$ws.Range("D44").Value = "'0.008530"
$ws.Range("E44").Value = "'-1.90%"

# Row 45
$ws.Range("D45").Value = "'0.00007091"
$ws.Range("E45").Value = "'-0.34%"

# Row 46
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.27%"

# Row 47
$ws.Range("E47").Value = "'54.99%"

# Row 48
$ws.Range("D48").Value = "'0.003110"
$ws.Range("E48").Value = "'-9.45%"

# Row 49
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.27%"

# Row 50
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.27%"
